# Regenerate the localization-status report for the "handoff" run:
#  - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#    timestamps forward (new handoff files were generated ~30s later).
#  - Mark the newly-handed-off rows with Priority "ht" (handoff type).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows 8,9,10,11,13,14 (row 12 / cd373b12-... is excluded - its handoff
# priority mismatched handoff type "ht", so it keeps its own timestamp).
$rows = @(8, 9, 10, 11, 13, 14)

foreach ($r in $rows) {
    # Overview!G<r> : "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-08-23 22:19:19"

    # de-de!H<r> : "Latest Handoff Datetime"
    $dede.Range("H$r").Value = "2016-08-23 22:19:19"

    # zh-cn!H<r> : "Latest Handoff Datetime"
    $zhcn.Range("H$r").Value = "2016-08-23 22:19:14"

    # Priority column (E) on both language sheets now reports "ht"
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
